$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: date 44309 -> 44305, volume 40 -> 50
$ws.Range("D5").Value = [DateTime]::FromOADate(44305)
$ws.Range("M5").Value = 50

# Row 6: date 44309 -> 44305, volume 70 -> 60
$ws.Range("D6").Value = [DateTime]::FromOADate(44305)
$ws.Range("M6").Value = 60

# Row 7: date 44305 -> 44309, volume 50 -> 40
$ws.Range("D7").Value = [DateTime]::FromOADate(44309)
$ws.Range("M7").Value = 40

# Row 8: date 44305 -> 44309, volume 60 -> 70
$ws.Range("D8").Value = [DateTime]::FromOADate(44309)
$ws.Range("M8").Value = 70
